$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Add sectors non waste" path in B4 to reflect the new repo location
$ws.Range("B4").Value = "C:\Users\loren\Documents\GitHub\Waste-MARIO\DWMRIO\waste-mario\Add sectors\Add_non_waste.xlsx"

# Move the active selection to B5 as reflected in the saved file
$ws.Range("B5").Select()
